# issue #5: stock data output to json file
#
# Changes applied to the "股票" (stock) worksheet:
#  1. Fix four company/owner-name typos that contained a stray embedded
#     space (OCR/typo artefacts) and normalise the full-width-comma
#     quantity value "20，000" -> "20000" (kept as text).
#  2. Insert a new "property_category" column right after "total" (i.e.
#     before the existing "date" column), shifting "date",
#     "legislator_name" and "legislator_id" one column to the right, and
#     fill the new column with the constant value "stock" for every data
#     row - this is the "stock data output to json file" change: each
#     property-type sheet gets tagged with its own category so the JSON
#     exporter downstream can tell rows apart once they're merged.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# --- 1. Typo / formatting fixes -------------------------------------------------
$ws.Cells.Item(2, 2).Value  = "太举洋電線電纜股份有限公司"
$ws.Cells.Item(7, 2).Value  = "春雨開發股份有限公司(原正華）"
$ws.Cells.Item(10, 2).Value = "太平洋電線電纜股份有限公司"
$ws.Cells.Item(13, 2).Value = "大成長城企業股份有限公司"

# Keep this one as TEXT ("20000"), matching the original cell's string type.
$qtyCell = $ws.Cells.Item(13, 4)
$qtyCell.NumberFormat = "@"
$qtyCell.Value = "20000"

# --- 2. Insert the new property_category column (H) -----------------------------
$ws.Columns.Item(8).Insert()

$ws.Cells.Item(1, 8).Value = "property_category"
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 8).Value = "stock"
}
